$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header for column F and copy header style/format from column E
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Timestamp values for rows 2-100 (column F), as text strings
$timeTaken = @(
  "2021-10-05 13:41:57.513861",
  "2021-10-05 13:41:57.513874",
  "2021-10-05 13:41:57.513877",
  "2021-10-05 13:41:57.513881",
  "2021-10-05 13:41:57.513884",
  "2021-10-05 13:41:57.513887",
  "2021-10-05 13:41:57.513890",
  "2021-10-05 13:41:57.513893",
  "2021-10-05 13:41:57.513896",
  "2021-10-05 13:41:57.513900",
  "2021-10-05 13:41:57.513903",
  "2021-10-05 13:41:57.513906",
  "2021-10-05 13:41:57.513909",
  "2021-10-05 13:41:57.513911",
  "2021-10-05 13:41:57.513914",
  "2021-10-05 13:41:57.513917",
  "2021-10-05 13:41:57.513921",
  "2021-10-05 13:41:57.513924",
  "2021-10-05 13:41:57.513927",
  "2021-10-05 13:41:57.513930",
  "2021-10-05 13:41:57.513933",
  "2021-10-05 13:41:57.513936",
  "2021-10-05 13:41:57.513939",
  "2021-10-05 13:41:57.513941",
  "2021-10-05 13:41:57.513945",
  "2021-10-05 13:41:57.513948",
  "2021-10-05 13:41:57.513951",
  "2021-10-05 13:41:57.513954",
  "2021-10-05 13:41:57.513957",
  "2021-10-05 13:41:57.513960",
  "2021-10-05 13:41:57.513963",
  "2021-10-05 13:41:57.513966",
  "2021-10-05 13:41:57.513969",
  "2021-10-05 13:41:57.513972",
  "2021-10-05 13:41:57.513975",
  "2021-10-05 13:41:57.513978",
  "2021-10-05 13:41:57.513981",
  "2021-10-05 13:41:57.513984",
  "2021-10-05 13:41:57.513987",
  "2021-10-05 13:41:57.513990",
  "2021-10-05 13:41:57.513993",
  "2021-10-05 13:41:57.513996",
  "2021-10-05 13:41:57.513999",
  "2021-10-05 13:41:57.514002",
  "2021-10-05 13:41:57.514005",
  "2021-10-05 13:41:57.514008",
  "2021-10-05 13:41:57.514011",
  "2021-10-05 13:41:57.514014",
  "2021-10-05 13:41:57.514017",
  "2021-10-05 13:41:57.514020",
  "2021-10-05 13:41:57.514023",
  "2021-10-05 13:41:57.514026",
  "2021-10-05 13:41:57.514029",
  "2021-10-05 13:41:57.514032",
  "2021-10-05 13:41:57.514036",
  "2021-10-05 13:41:57.514039",
  "2021-10-05 13:41:57.514042",
  "2021-10-05 13:41:57.514045",
  "2021-10-05 13:41:57.514047",
  "2021-10-05 13:41:57.514050",
  "2021-10-05 13:41:57.514053",
  "2021-10-05 13:41:57.514056",
  "2021-10-05 13:41:57.514059",
  "2021-10-05 13:41:57.514062",
  "2021-10-05 13:41:57.514067",
  "2021-10-05 13:41:57.514070",
  "2021-10-05 13:41:57.514073",
  "2021-10-05 13:41:57.514076",
  "2021-10-05 13:41:57.514079",
  "2021-10-05 13:41:57.514082",
  "2021-10-05 13:41:57.514085",
  "2021-10-05 13:41:57.514088",
  "2021-10-05 13:41:57.514091",
  "2021-10-05 13:41:57.514094",
  "2021-10-05 13:41:57.514097",
  "2021-10-05 13:41:57.514100",
  "2021-10-05 13:41:57.514106",
  "2021-10-05 13:41:57.514109",
  "2021-10-05 13:41:57.514112",
  "2021-10-05 13:41:57.514115",
  "2021-10-05 13:41:57.514118",
  "2021-10-05 13:41:57.514122",
  "2021-10-05 13:41:57.514125",
  "2021-10-05 13:41:57.514128",
  "2021-10-05 13:41:57.514131",
  "2021-10-05 13:41:57.514134",
  "2021-10-05 13:41:57.514137",
  "2021-10-05 13:41:57.514140",
  "2021-10-05 13:41:57.514143",
  "2021-10-05 13:41:57.514146",
  "2021-10-05 13:41:57.514149",
  "2021-10-05 13:41:57.514152",
  "2021-10-05 13:41:57.514156",
  "2021-10-05 13:41:57.514160",
  "2021-10-05 13:41:57.514163",
  "2021-10-05 13:41:57.514166",
  "2021-10-05 13:41:57.514169",
  "2021-10-05 13:41:57.514172",
  "2021-10-05 13:41:57.514175"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

$wb.Save()